# Weekly price-sheet update: a new weekly record is inserted as row 186,
# pushing the existing rows 186..269 down to 187..270 (the last existing
# row, formerly 269, becomes row 270).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 186; Excel shifts rows 186-269 down to 187-270
# and carries their values/formatting with them automatically.
$ws.Rows.Item(186).Insert()

# Populate the newly inserted row 186 with the new weekly record.
$ws.Range("A186").Value = 7
$ws.Range("B186").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C186").Value = "Ñuble"
$ws.Range("D186").Value = 44636
$ws.Range("E186").Value = 16
$ws.Range("F186").Value = 100114013
$ws.Range("G186").Value = "Zanahoria"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 120
$ws.Range("K186").Value = 6500
$ws.Range("L186").Value = 7000
$ws.Range("M186").Value = 6750
$ws.Range("N186").Value = "`$/saco 20 kilos"
$ws.Range("O186").Value = "Provincia de Diguillín"
$ws.Range("P186").Value = 338
$ws.Range("Q186").Value = 20
$ws.Range("R186").Value = "Hortaliza"
